$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the two new columns
$ws.Range("D1").Value = "color1"
$ws.Range("E1").Value = "color2"
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("D1:E1").VerticalAlignment = -4108

# Per-team color1 / color2 hex codes (data rows 2-45)
$ws.Range("D2").Value = "#007242"
$ws.Range("E2").Value = "#ffffff"
$ws.Range("D3").Value = "#ff1100"
$ws.Range("E3").Value = "#ffffff"
$ws.Range("D4").Value = "#e43031"
$ws.Range("E4").Value = "#231f20"
$ws.Range("D5").Value = "#dc1212"
$ws.Range("E5").Value = "#010101"
$ws.Range("D6").Value = "#1f1a17"
$ws.Range("E6").Value = "#ffffff"
$ws.Range("D7").Value = "#0b5b8d"
$ws.Range("E7").Value = "#ffffff"
$ws.Range("D8").Value = "#2762a6"
$ws.Range("E8").Value = "#cd3529"
$ws.Range("D9").Value = "#f5e145"
$ws.Range("E9").Value = "#cf192e"
$ws.Range("D10").Value = "#000000"
$ws.Range("E10").Value = "#ffffff"
$ws.Range("D11").Value = "#0f2d5b"
$ws.Range("E11").Value = "#e41349"
$ws.Range("D12").Value = "#f6dc01"
$ws.Range("E12").Value = "#009b45"
$ws.Range("D13").Value = "#000000"
$ws.Range("E13").Value = "#fefefe"
$ws.Range("D14").Value = "#026b38"
$ws.Range("E14").Value = "#ffffff"
$ws.Range("D15").Value = "#1e1214"
$ws.Range("E15").Value = "#d82531"
$ws.Range("D16").Value = "#006755"
$ws.Range("E16").Value = "#ffffff"
$ws.Range("D17").Value = "#f4cf16"
$ws.Range("E17").Value = "#1f150a"
$ws.Range("D18").Value = "#1e3d8e"
$ws.Range("E18").Value = "#ffffff"
$ws.Range("D19").Value = "#041e5c"
$ws.Range("E19").Value = "#ffffff"
$ws.Range("D20").Value = "#016c32"
$ws.Range("E20").Value = "#ffdb01"
$ws.Range("D21").Value = "#231f20"
$ws.Range("E21").Value = "#03a550"
$ws.Range("D22").Value = "#c3281e"
$ws.Range("E22").Value = "#070308"
$ws.Range("D23").Value = "#870a28"
$ws.Range("E23").Value = "#00613c"
$ws.Range("D24").Value = "#1f5ea1"
$ws.Range("E24").Value = "#e1251b"
$ws.Range("D25").Value = "#26603c"
$ws.Range("E25").Value = "#ffffff"
$ws.Range("D26").Value = "#0e93d2"
$ws.Range("E26").Value = "#231f20"
$ws.Range("D27").Value = "#006e4d"
$ws.Range("E27").Value = "#ffffff"
$ws.Range("D28").Value = "#e5050f"
$ws.Range("E28").Value = "#ffffff"
$ws.Range("D29").Value = "#64b054"
$ws.Range("E29").Value = "#cc262d"
$ws.Range("D30").Value = "#ee3338"
$ws.Range("E30").Value = "#ffffff"
$ws.Range("D31").Value = "#118e4c"
$ws.Range("E31").Value = "#ffffff"
$ws.Range("D32").Value = "#ed2129"
$ws.Range("E32").Value = "#ffffff"
$ws.Range("D33").Value = "#006437"
$ws.Range("E33").Value = "#ffffff"
$ws.Range("D34").Value = "#cd3529"
$ws.Range("E34").Value = "#304296"
$ws.Range("D35").Value = "#263484"
$ws.Range("E35").Value = "#fbfbfb"
$ws.Range("D36").Value = "#000000"
$ws.Range("E36").Value = "#fdfdfd"
$ws.Range("D37").Value = "#198538"
$ws.Range("E37").Value = "#ee2722"
$ws.Range("D38").Value = "#231f20"
$ws.Range("E38").Value = "#ee3124"
$ws.Range("D39").Value = "#267946"
$ws.Range("E39").Value = "#2d3492"
$ws.Range("D40").Value = "#000000"
$ws.Range("E40").Value = "#ffffff"
$ws.Range("D41").Value = "#034ea2"
$ws.Range("E41").Value = "#00a651"
$ws.Range("D42").Value = "#000000"
$ws.Range("E42").Value = "#ed1c24"
$ws.Range("D43").Value = "#d91a21"
$ws.Range("E43").Value = "#ffda00"
$ws.Range("D44").Value = "#000000"
$ws.Range("E44").Value = "#ffffff"
$ws.Range("D45").Value = "#ff1100"
$ws.Range("E45").Value = "#000000"

# Restore the active selection as recorded after the edit
$ws.Range("C6").Select() | Out-Null
